$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 100007096
$ws.Cells.Item(64, 9).Value = 8330.666999999999
$ws.Cells.Item(64, 10).Value = 142863710
$ws.Cells.Item(64, 11).Value = 8330.666999999999
$ws.Cells.Item(64, 12).Value = 142863710
$ws.Cells.Item(64, 13).Value = -8082.666999999999
$ws.Cells.Item(64, 14).Value = -142864206

$ws.Cells.Item(67, 8).Value = 100007096
$ws.Cells.Item(67, 9).Value = 8330.666999999999
$ws.Cells.Item(67, 10).Value = 142863710
$ws.Cells.Item(67, 11).Value = 8330.666999999999
$ws.Cells.Item(67, 12).Value = 142863710
$ws.Cells.Item(67, 13).Value = -7472.666999999999
$ws.Cells.Item(67, 14).Value = -142865426

$ws.Cells.Item(107, 8).Value = 289.3684
$ws.Cells.Item(107, 9).Value = 227.33333
$ws.Cells.Item(107, 10).Value = 1406
$ws.Cells.Item(107, 11).Value = 227.33333
$ws.Cells.Item(107, 12).Value = 1406
$ws.Cells.Item(107, 13).Value = 1692.66667

$ws.Cells.Item(112, 8).Value = 1992.2667
$ws.Cells.Item(112, 9).Value = 1366.6666
$ws.Cells.Item(112, 10).Value = 2061.7778
$ws.Cells.Item(112, 11).Value = 4099.9998
$ws.Cells.Item(112, 12).Value = 6185.3334
$ws.Cells.Item(112, 13).Value = -2991.9998
$ws.Cells.Item(112, 14).Value = -8401.3334

$ws.Cells.Item(132, 8).Value = 6180.4517
$ws.Cells.Item(132, 9).Value = 6859.92
$ws.Cells.Item(132, 10).Value = 3349.3333
$ws.Cells.Item(132, 11).Value = 20579.76
$ws.Cells.Item(132, 12).Value = 10047.9999
$ws.Cells.Item(132, 13).Value = -18049.76
$ws.Cells.Item(132, 14).Value = -15107.9999

$ws.Cells.Item(137, 8).Value = 1256011.6
$ws.Cells.Item(137, 9).Value = 4167692
$ws.Cells.Item(137, 10).Value = 8148.5
$ws.Cells.Item(137, 11).Value = 12503076
$ws.Cells.Item(137, 12).Value = 24445.5
$ws.Cells.Item(137, 13).Value = -12500526
$ws.Cells.Item(137, 14).Value = -29545.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 464.66666
$ws.Cells.Item(2, 9).Value = 359.7619
$ws.Cells.Item(2, 10).Value = 1199
$ws.Cells.Item(2, 11).Value = 359.7619
$ws.Cells.Item(2, 12).Value = 1199
$ws.Cells.Item(2, 13).Value = -246.7619

$ws.Cells.Item(32, 8).Value = 3729.8518
$ws.Cells.Item(32, 9).Value = 3123.182
$ws.Cells.Item(32, 10).Value = 6399.2
$ws.Cells.Item(32, 11).Value = 3123.182
$ws.Cells.Item(32, 12).Value = 6399.2
$ws.Cells.Item(32, 13).Value = -2836.182
$ws.Cells.Item(32, 14).Value = -6973.2

$ws.Cells.Item(63, 8).Value = 3499.5
$ws.Cells.Item(63, 9).Value = 3000
$ws.Cells.Item(63, 10).Value = 3999
$ws.Cells.Item(63, 11).Value = 3000
$ws.Cells.Item(63, 12).Value = 3999
$ws.Cells.Item(63, 13).Value = -2314

$ws.Cells.Item(66, 8).Value = 3499.5
$ws.Cells.Item(66, 9).Value = 3000
$ws.Cells.Item(66, 10).Value = 3999
$ws.Cells.Item(66, 11).Value = 15000
$ws.Cells.Item(66, 12).Value = 19995
$ws.Cells.Item(66, 13).Value = -11568

$ws.Cells.Item(116, 8).Value = 464.66666
$ws.Cells.Item(116, 9).Value = 359.7619
$ws.Cells.Item(116, 10).Value = 1199
$ws.Cells.Item(116, 11).Value = 359.7619
$ws.Cells.Item(116, 12).Value = 1199
$ws.Cells.Item(116, 13).Value = 1934.2381

$ws.Cells.Item(132, 8).Value = 1738.1351
$ws.Cells.Item(132, 9).Value = 1139.0968
$ws.Cells.Item(132, 10).Value = 4833.1665
$ws.Cells.Item(132, 11).Value = 3417.2904
$ws.Cells.Item(132, 12).Value = 14499.4995
$ws.Cells.Item(132, 13).Value = -887.2903999999999

$ws.Cells.Item(135, 8).Value = 110999
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 110999
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 110999
$ws.Cells.Item(135, 14).Value = -121139

$ws.Cells.Item(139, 8).Value = 82330.836
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 82330.836
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 82330.836
$ws.Cells.Item(139, 14).Value = -92610.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 464.66666
$ws.Cells.Item(3, 9).Value = 359.7619
$ws.Cells.Item(3, 10).Value = 1199
$ws.Cells.Item(3, 11).Value = 359.7619
$ws.Cells.Item(3, 12).Value = 1199
$ws.Cells.Item(3, 13).Value = -245.7619

$ws.Cells.Item(134, 8).Value = 4125.1665
$ws.Cells.Item(134, 9).Value = 4076.3044
$ws.Cells.Item(134, 10).Value = 4285.7144
$ws.Cells.Item(134, 11).Value = 12228.9132
$ws.Cells.Item(134, 12).Value = 12857.1432
$ws.Cells.Item(134, 13).Value = -9693.913199999999

$ws.Cells.Item(139, 8).Value = 52129.5
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 52129.5
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 52129.5
$ws.Cells.Item(139, 14).Value = -62409.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4374.735
$ws.Cells.Item(31, 9).Value = 2979.348
$ws.Cells.Item(31, 10).Value = 5609.115
$ws.Cells.Item(31, 11).Value = 2979.348
$ws.Cells.Item(31, 12).Value = 5609.115
$ws.Cells.Item(31, 13).Value = -2684.348
$ws.Cells.Item(31, 14).Value = -6199.115

$ws.Cells.Item(34, 8).Value = 4374.735
$ws.Cells.Item(34, 9).Value = 2979.348
$ws.Cells.Item(34, 10).Value = 5609.115
$ws.Cells.Item(34, 11).Value = 2979.348
$ws.Cells.Item(34, 12).Value = 5609.115
$ws.Cells.Item(34, 13).Value = -2777.348
$ws.Cells.Item(34, 14).Value = -6013.115

$ws.Cells.Item(62, 8).Value = 12507685
$ws.Cells.Item(62, 9).Value = 25008876
$ws.Cells.Item(62, 10).Value = 6495.75
$ws.Cells.Item(62, 11).Value = 25008876
$ws.Cells.Item(62, 12).Value = 6495.75
$ws.Cells.Item(62, 13).Value = -25008252
$ws.Cells.Item(62, 14).Value = -7743.75

$ws.Cells.Item(65, 8).Value = 12507685
$ws.Cells.Item(65, 9).Value = 25008876
$ws.Cells.Item(65, 10).Value = 6495.75
$ws.Cells.Item(65, 11).Value = 125044380
$ws.Cells.Item(65, 12).Value = 32478.75
$ws.Cells.Item(65, 13).Value = -125041260
$ws.Cells.Item(65, 14).Value = -38718.75

$ws.Cells.Item(74, 8).Value = 67157
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 67157
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 67157
$ws.Cells.Item(74, 14).Value = -68905

$ws.Cells.Item(77, 8).Value = 67157
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 67157
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 201471
$ws.Cells.Item(77, 14).Value = -210207

$ws.Cells.Item(82, 8).Value = 40000
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 40000
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 40000
$ws.Cells.Item(82, 14).Value = -40722

$ws.Cells.Item(85, 8).Value = 40000
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 40000
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 40000
$ws.Cells.Item(85, 14).Value = -42496

$ws.Cells.Item(99, 8).Value = 4399.875
$ws.Cells.Item(99, 9).Value = 4280
$ws.Cells.Item(99, 10).Value = 4599.6665
$ws.Cells.Item(99, 11).Value = 4280
$ws.Cells.Item(99, 12).Value = 4599.6665
$ws.Cells.Item(99, 13).Value = -2782
$ws.Cells.Item(99, 14).Value = -7595.6665

$ws.Cells.Item(106, 8).Value = 65037
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 65037
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 65037
$ws.Cells.Item(106, 14).Value = -67561

$ws.Cells.Item(122, 8).Value = 4637.2173
$ws.Cells.Item(122, 9).Value = 3825.0833
$ws.Cells.Item(122, 10).Value = 5523.1816
$ws.Cells.Item(122, 11).Value = 11475.2499
$ws.Cells.Item(122, 12).Value = 16569.5448
$ws.Cells.Item(122, 13).Value = -9025.249899999999

$ws.Cells.Item(126, 8).Value = 4399.875
$ws.Cells.Item(126, 9).Value = 4280
$ws.Cells.Item(126, 10).Value = 4599.6665
$ws.Cells.Item(126, 11).Value = 12840
$ws.Cells.Item(126, 12).Value = 13798.9995
$ws.Cells.Item(126, 13).Value = -10370
$ws.Cells.Item(126, 14).Value = -18738.9995

$ws.Cells.Item(135, 8).Value = 121998.2
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 121998.2
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 121998.2
$ws.Cells.Item(135, 14).Value = -132138.2

$ws.Cells.Item(138, 8).Value = 68397.5
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 68397.5
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 68397.5
$ws.Cells.Item(138, 14).Value = -78677.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 1056040.9
$ws.Cells.Item(60, 9).Value = 3333749.5
$ws.Cells.Item(60, 10).Value = 4790.769
$ws.Cells.Item(60, 11).Value = 10001248.5
$ws.Cells.Item(60, 12).Value = 14372.307
$ws.Cells.Item(60, 13).Value = -10000997.5
$ws.Cells.Item(60, 14).Value = -14874.307

$ws.Cells.Item(80, 8).Value = 3307.25
$ws.Cells.Item(80, 9).Value = 3474
$ws.Cells.Item(80, 10).Value = 3251.6667
$ws.Cells.Item(80, 11).Value = 10422
$ws.Cells.Item(80, 12).Value = 9755.000100000001
$ws.Cells.Item(80, 13).Value = -9486
$ws.Cells.Item(80, 14).Value = -11627.0001

$ws.Cells.Item(83, 8).Value = 3307.25
$ws.Cells.Item(83, 9).Value = 3474
$ws.Cells.Item(83, 10).Value = 3251.6667
$ws.Cells.Item(83, 11).Value = 31266
$ws.Cells.Item(83, 12).Value = 29265.0003
$ws.Cells.Item(83, 13).Value = -26586
$ws.Cells.Item(83, 14).Value = -38625.0003

$ws.Cells.Item(122, 8).Value = 629.2174
$ws.Cells.Item(122, 9).Value = 919.1667
$ws.Cells.Item(122, 10).Value = 526.8823
$ws.Cells.Item(122, 11).Value = 8272.5003
$ws.Cells.Item(122, 12).Value = 4741.9407
$ws.Cells.Item(122, 13).Value = -5822.5003
$ws.Cells.Item(122, 14).Value = -9641.940699999999

$ws.Cells.Item(132, 8).Value = 2071.9412
$ws.Cells.Item(132, 9).Value = 1258.7778
$ws.Cells.Item(132, 10).Value = 2986.75
$ws.Cells.Item(132, 11).Value = 11329.0002
$ws.Cells.Item(132, 12).Value = 26880.75
$ws.Cells.Item(132, 13).Value = -8799.0002
$ws.Cells.Item(132, 14).Value = -31940.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 10874.5625
$ws.Cells.Item(24, 9).Value = 10996.5
$ws.Cells.Item(24, 10).Value = 10857.143
$ws.Cells.Item(24, 11).Value = 10996.5
$ws.Cells.Item(24, 12).Value = 10857.143
$ws.Cells.Item(24, 13).Value = -10823.5
$ws.Cells.Item(24, 14).Value = -11203.143

$ws.Cells.Item(132, 8).Value = 3703.1538
$ws.Cells.Item(132, 9).Value = 3377.4285
$ws.Cells.Item(132, 10).Value = 4083.1667
$ws.Cells.Item(132, 11).Value = 10132.2855
$ws.Cells.Item(132, 12).Value = 12249.5001
$ws.Cells.Item(132, 13).Value = -7602.2855
$ws.Cells.Item(132, 14).Value = -17309.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 14999.5
$ws.Cells.Item(14, 9).Value = 14999
$ws.Cells.Item(14, 10).Value = 15000
$ws.Cells.Item(14, 11).Value = 14999
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = -14827
$ws.Cells.Item(14, 14).Value = -15344

$ws.Cells.Item(23, 8).Value = 10000
$ws.Cells.Item(23, 9).Value = 10000
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 10000
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -9770
$ws.Cells.Item(23, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 3850
$ws.Cells.Item(46, 9).Value = 5750
$ws.Cells.Item(46, 10).Value = 1950
$ws.Cells.Item(46, 11).Value = 5750
$ws.Cells.Item(46, 12).Value = 1950
$ws.Cells.Item(46, 13).Value = -5562
$ws.Cells.Item(46, 14).Value = -2326

$ws.Cells.Item(100, 8).Value = 1703.6666
$ws.Cells.Item(100, 9).Value = 1399
$ws.Cells.Item(100, 10).Value = 2313
$ws.Cells.Item(100, 11).Value = 1399
$ws.Cells.Item(100, 12).Value = 2313
$ws.Cells.Item(100, 13).Value = -858

$ws.Cells.Item(122, 8).Value = 2770.8572
$ws.Cells.Item(122, 9).Value = 2791.6924
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 8375.0772
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -5925.0772

$ws.Cells.Item(136, 8).Value = 5395.3
$ws.Cells.Item(136, 9).Value = 6325.6665
$ws.Cells.Item(136, 10).Value = 3999.75
$ws.Cells.Item(136, 11).Value = 18976.9995
$ws.Cells.Item(136, 12).Value = 11999.25
$ws.Cells.Item(136, 13).Value = -16426.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 65228.125
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 65228.125
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 65228.125
$ws.Cells.Item(46, 14).Value = -65690.125

$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(101, 8).Value = 59624.5
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 59624.5
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 59624.5
$ws.Cells.Item(101, 14).Value = -66114.5

$ws.Cells.Item(104, 8).Value = 19745.5
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 19745.5
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 19745.5
$ws.Cells.Item(104, 14).Value = -26733.5

$ws.Cells.Item(134, 8).Value = 65228.125
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 65228.125
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 195684.375
$ws.Cells.Item(134, 14).Value = -200754.375

$ws.Cells.Item(136, 8).Value = 76927560
$ws.Cells.Item(136, 9).Value = 142857800
$ws.Cells.Item(136, 10).Value = 8966.666999999999
$ws.Cells.Item(136, 11).Value = 428573400
$ws.Cells.Item(136, 12).Value = 26900.001
$ws.Cells.Item(136, 13).Value = -428570850
